$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("M:M").Delete() | Out-Null
$ws.Range("M1").Select() | Out-Null
